# Companies.xlsx update — add eight new companies, reorder/relabel the
# header columns, restyle the header (bold 12pt) and the data cells
# (wrap text, mmmm-yyyy date format), switch two rows to an Arial 10pt
# company-name font, resize columns and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header row (row 1): re-label / re-order columns B..H and make the
#    whole header bold 12pt (A1 keeps no wrap, B1:H1 keep their wrap).
# ---------------------------------------------------------------------
$ws.Range("A1").Value  = "Company Name"
$ws.Range("B1").Value  = "scode"
$ws.Range("C1").Value  = "scriptcd"
$ws.Range("D1").Value  = "comName"
$ws.Range("E1").Value  = "sname1"
$ws.Range("F1").Value  = "sname2"
$ws.Range("G1").Value  = "Qname"
$ws.Range("H1").Value  = "qtrid"

$ws.Range("A1:H1").Font.Size = 12
$ws.Range("A1:H1").Font.Bold = $true

# ---------------------------------------------------------------------
# 2. Existing data row (row 2): reorder scode/scriptcd and the
#    sname1/sname2/Qname/qtrid block to match the new header, then drop
#    the old 10pt font back to the sheet's normal 12pt font so the row
#    keeps only the wrap-text / date formatting.
# ---------------------------------------------------------------------
$ws.Range("B2").Value = 500387
$ws.Range("C2").Value = 324875
$ws.Range("E2").Value = "shree-cement-ltd"
$ws.Range("F2").Value = "shreecem"
$ws.Range("G2").Value = 45627
$ws.Range("H2").Value = 124

$ws.Range("A2:H2").Font.Size = 12

# ---------------------------------------------------------------------
# 3. New rows 3-10: eight additional companies.
# ---------------------------------------------------------------------
$data = @(
    @("Tata Motors",        500570, 327610, "TATA MOTORS LTD.",          "tata-motors-ltd",           "tatamotors", 124),
    @("PC Jeweller",        534809, 327230, "PC JEWELLER LTD.",          "pc-jeweller-ltd",           "pcjeweller", 124),
    @("Vakrangee",          511431, 325794, "Vakrangee Limited-$",       "vakrangee-limited",         "vakrangee",  124),
    @("Jubilant Foodworks", 533155, 326075, "Jubilant FoodWorks Ltd",    "jubilant-foodworks-ltd",    "jublfood",   124),
    @("Paytm",              543396, 325689, "One 97 Communications Ltd","one-97-communications-ltd", "paytm",      124),
    @("Sanghvi Movers",     530073, 327162, "SANGHVI MOVERS LTD.",       "sanghvi-movers-ltd",        "sanghvimov", 124),
    @("MGL",                539957, 326504, "Mahanagar Gas Ltd",         "mahanagar-gas-ltd",         "mgl",        124),
    @("Religare",           532915, 328115, "RELIGARE ENTERPRISES LTD.","religare-enterprises-ltd",  "religare",   124)
)

$r = 3
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = 45627
    $ws.Cells.Item($r, 8).Value = $row[6]
    $ws.Rows.Item($r).RowHeight = 15.75
    $r = $r + 1
}

# Row 2 also gets the "new row" height used throughout the table.
$ws.Rows.Item(1).RowHeight = 17
$ws.Rows.Item(2).RowHeight = 17

# ---------------------------------------------------------------------
# 4. Formatting: row 2 (pre-existing row) keeps its wrap text on every
#    value column; the freshly added rows 3-10 only carry it on the
#    qtrid column, plus every data row gets the mmmm-yyyy date format
#    (with wrap) on the date column.
# ---------------------------------------------------------------------
$ws.Range("B2:F2").WrapText = $true
$ws.Range("H2:H10").WrapText = $true
$ws.Range("G2:G10").WrapText = $true
$ws.Range("G2:G10").NumberFormat = "mmmm\ yyyy"

# ---------------------------------------------------------------------
# 5. Company-name font override for MGL (row 9) and Religare (row 10):
#    Arial 10pt, matching the rest of the old "small" font family.
# ---------------------------------------------------------------------
$ws.Range("A9:A10").Font.Size = 10
$ws.Range("A9:A10").Font.Name = "Arial"

# ---------------------------------------------------------------------
# 6. Column widths (best-fit in the source workbook).
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 16.1666666666667
$ws.Columns.Item(2).ColumnWidth = 6.33333333333333
$ws.Columns.Item(3).ColumnWidth = 6.66666666666667
$ws.Columns.Item(4).ColumnWidth = 19.5
$ws.Columns.Item(5).ColumnWidth = 18.6666666666667
$ws.Columns.Item(6).ColumnWidth = 9.5
$ws.Columns.Item(7).ColumnWidth = 13.1666666666667
$ws.Columns.Item(8).ColumnWidth = 4.16666666666667

# ---------------------------------------------------------------------
# 7. Selection used by the author when they last saved the file.
# ---------------------------------------------------------------------
[void]$ws.Range("K9").Select()
